$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.854.71"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").Value = "2.566.92"
$ws.Range("E3").Value = "  +1.93%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.54"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.88"
$ws.Range("E6").Value = "  +7.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.571"
$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  +0.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.42"
$ws.Range("E10").Value = "  +3.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.44"
$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.109"
$ws.Range("E13").Value = "  +0.72%  "

$ws.Range("D14").Value = "2.955.06"
$ws.Range("E14").Value = "  +2.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.11"
$ws.Range("E15").Value = "  +8.20%  "

$ws.Range("D16").Value = "2.539.02"
$ws.Range("E16").Value = "  +2.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.849"
$ws.Range("E17").Value = "  +1.82%  "

$ws.Range("D18").Value = "42.859.87"
$ws.Range("E18").Value = "  +1.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.85"
$ws.Range("E19").Value = "  +1.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.46"
$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("D21").Value = "0.0₃0957"
$ws.Range("E21").Value = "  +0.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.34"
$ws.Range("E22").Value = "  +0.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.68"
$ws.Range("E23").Value = "  -2.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.94"
$ws.Range("E24").Value = "  +0.94%  "

$ws.Range("E25").Value = "  +2.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.59"
$ws.Range("E26").Value = "  +0.75%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.11"
$ws.Range("E28").Value = "  +3.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.37"
$ws.Range("E29").Value = "  -1.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.16"
$ws.Range("E30").Value = "  +0.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.26"
$ws.Range("E31").Value = "  +1.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.73"
$ws.Range("E32").Value = "  -0.97%  "

$ws.Range("E33").Value = "  +14.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0807"
$ws.Range("E34").Value = "  +2.61%  "

$ws.Range("E35").Value = "  +0.91%  "

$ws.Range("E36").Value = "  -2.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.22"
$ws.Range("E37").Value = "  -1.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.40"
$ws.Range("E38").Value = "  -4.66%  "

$ws.Range("E39").Value = "  +0.34%  "

$ws.Range("E40").Value = "  +0.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.26"
$ws.Range("E41").Value = "  +13.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.31"
$ws.Range("E42").Value = "  +0.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.34"
$ws.Range("E43").Value = "  +4.16%  "

$ws.Range("E44").Value = "  +0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0299"
$ws.Range("E45").Value = "  +0.29%  "

$ws.Range("D46").Value = "1.974.17"
$ws.Range("E46").Value = "  -0.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.94"
$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("D48").Value = "2.806.61"
$ws.Range("E48").Value = "  +2.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "81.56"
$ws.Range("E49").Value = "  -2.53%  "

$ws.Range("E50").Value = "  +2.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.77"
$ws.Range("E51").Value = "  +0.35%  "
